# Apply cryptocurrency price/volume updates to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to remain a plain text string even when the
    # value looks numeric (e.g. "215.32"), then restore the default
    # (unstyled) cell style so no stray formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "25.910.04"
$ws.Range("E2").Value = "  +0.40%  "

Set-TextValue $ws.Range("D3") "1.640.79"
$ws.Range("E3").Value = "  +0.44%  "

Set-TextValue $ws.Range("D4") "1.004"
$ws.Range("E4").Value = "  +0.08%  "

Set-TextValue $ws.Range("D5") "215.32"
$ws.Range("E5").Value = "  +0.63%  "

Set-TextValue $ws.Range("D6") "0.5085"
$ws.Range("E6").Value = "  +1.33%  "

Set-TextValue $ws.Range("D7") "1.004"
$ws.Range("E7").Value = "  +0.34%  "

Set-TextValue $ws.Range("D8") "0.2566"
$ws.Range("E8").Value = "  +0.29%  "

Set-TextValue $ws.Range("D9") "0.06388"
$ws.Range("E9").Value = "  +0.33%  "

Set-TextValue $ws.Range("D10") "19.53"
$ws.Range("E10").Value = "  -0.54%  "

Set-TextValue $ws.Range("D11") "0.07756"
$ws.Range("E11").Value = "  +0.56%  "

Set-TextValue $ws.Range("D12") "4.306"
$ws.Range("E12").Value = "  +1.20%  "

Set-TextValue $ws.Range("D13") "1.649.27"
$ws.Range("E13").Value = "  +0.80%  "

Set-TextValue $ws.Range("D14") "0.5453"
$ws.Range("E14").Value = "  +0.58%  "

Set-TextValue $ws.Range("D15") "0.0₅7849"
$ws.Range("E15").Value = "  -0.29%  "

Set-TextValue $ws.Range("D16") "64.68"
$ws.Range("E16").Value = "  +1.81%  "

Set-TextValue $ws.Range("D17") "25.968.34"
$ws.Range("E17").Value = "  +0.61%  "

Set-TextValue $ws.Range("D18") "1.004"
$ws.Range("E18").Value = "  +0.11%  "

Set-TextValue $ws.Range("D19") "197.67"
$ws.Range("E19").Value = "  -1.28%  "

Set-TextValue $ws.Range("D20") "4.432"
$ws.Range("E20").Value = "  +2.11%  "

Set-TextValue $ws.Range("D21") "9.953"
$ws.Range("E21").Value = "  +0.69%  "

Set-TextValue $ws.Range("D22") "6.039"
$ws.Range("E22").Value = "  +1.82%  "

Set-TextValue $ws.Range("D23") "1.006"
$ws.Range("E23").Value = "  +0.42%  "

Set-TextValue $ws.Range("D24") "1.872"
$ws.Range("E24").Value = "  -2.51%  "

Set-TextValue $ws.Range("D25") "140.62"
$ws.Range("E25").Value = "  -0.15%  "

Set-TextValue $ws.Range("D26") "0.1145"
$ws.Range("E26").Value = "  +1.07%  "

Set-TextValue $ws.Range("D27") "6.887"
$ws.Range("E27").Value = "  +3.14%  "

Set-TextValue $ws.Range("D28") "15.71"

Set-TextValue $ws.Range("D29") "1.236"
$ws.Range("E29").Value = "  -0.02%  "

Set-TextValue $ws.Range("D30") "0.05021"
$ws.Range("E30").Value = "  +0.96%  "

Set-TextValue $ws.Range("D31") "3.260"
$ws.Range("E31").Value = "  -0.05%  "

Set-TextValue $ws.Range("D32") "3.180"
$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("E33").Value = "  +0.24%  "

Set-TextValue $ws.Range("D34") "2.361"
$ws.Range("E34").Value = "  -0.45%  "

Set-TextValue $ws.Range("D35") "0.8926"
$ws.Range("E35").Value = "  +0.80%  "

Set-TextValue $ws.Range("D36") "2.584"
$ws.Range("E36").Value = "  -1.67%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "0.5502"
$ws.Range("E37").Value = "  -0.87%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D38") "1.125.47"
$ws.Range("E38").Value = "  -3.38%  "

Set-TextValue $ws.Range("D39") "0.01557"
$ws.Range("E39").Value = "  +0.14%  "

Set-TextValue $ws.Range("D40") "1.005"
$ws.Range("E40").Value = "  +0.53%  "

Set-TextValue $ws.Range("D41") "2.546"
$ws.Range("E41").Value = "  -0.29%  "

Set-TextValue $ws.Range("D42") "0.0₈130"
$ws.Range("E42").Value = "  +12.37%  "

Set-TextValue $ws.Range("D43") "5.628"
$ws.Range("E43").Value = "  -0.65%  "

Set-TextValue $ws.Range("D44") "0.8166"
$ws.Range("E44").Value = "  +1.63%  "

Set-TextValue $ws.Range("D45") "99.94"
$ws.Range("E45").Value = "  +0.68%  "

Set-TextValue $ws.Range("D46") "1.778.84"
$ws.Range("E46").Value = "  +0.51%  "

Set-TextValue $ws.Range("D47") "0.4526"
$ws.Range("E47").Value = "  +0.13%  "

Set-TextValue $ws.Range("D48") "0.9986"
$ws.Range("E48").Value = "  +0.26%  "

Set-TextValue $ws.Range("D49") "54.77"
$ws.Range("E49").Value = "  +0.08%  "

Set-TextValue $ws.Range("D50") "0.05083"
$ws.Range("E50").Value = "  +0.29%  "

Set-TextValue $ws.Range("D51") "1.005"
$ws.Range("E51").Value = "  +0.39%  "
